# Applies the "additional scraping" update described in the commit:
#   1. Renames the "ODI Batting" column D header from MATCH_CARD_LINK to
#      MATCH_CODE and shrinks each row's full scorecard URL down to just the
#      numeric match code that was embedded in it.
#   2. Adds a brand-new "Player Info" worksheet ahead of "ODI Batting" that
#      carries the player's static metadata (ID, NAME, BATTING_HAND,
#      BOWL_STYLE).

$wb = $excel.ActiveWorkbook

# --- 1. Existing "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$batting = $wb.Worksheets.Item("ODI Batting")

$batting.Range("D1").Value = "MATCH_CODE"

$matchCodes = @("4189", "4190", "4257", "4259", "4262", "4265", "4290", "4299", "4306", "4309", "4315", "4323", "4332", "4377", "4378", "4379")

# Format the column as text first so these numeric-looking codes are kept as
# text (same as every other "numeric" column in this sheet) instead of being
# auto-converted to a Number.
$batting.Range("D2:D17").NumberFormat = "@"
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $batting.Cells.Item($row, 4).Value = $matchCodes[$i]
}

# --- 2. Add the new "Player Info" sheet ahead of "ODI Batting" ---
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Header row, styled like the header row on "ODI Batting" (bold, bordered,
# centered horizontally, top-aligned vertically).
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4615"
$playerInfo.Range("B2").Value = "Hazratullah Zazai"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# Match the page margins used throughout the rest of the workbook
# (0.75in/0.75in/1in/1in, 0.5in header/footer == 54/54/72/72/36/36 points).
$playerInfo.PageSetup.LeftMargin = 54
$playerInfo.PageSetup.RightMargin = 54
$playerInfo.PageSetup.TopMargin = 72
$playerInfo.PageSetup.BottomMargin = 72
$playerInfo.PageSetup.HeaderMargin = 36
$playerInfo.PageSetup.FooterMargin = 36

$playerInfo.Move($batting)
